$d = $word.ActiveDocument

# 1. "volume that are in the template" -> "volumes that are in the template"
$d.Content.Find.Execute("volume that are in the template", $true, $false, $false, $false, $false, `
    $true, 1, $false, "volumes that are in the template", 2) | Out-Null

# 2. Remove "HAS TO BE DONE MANUALLY" at end of paragraph 1 (keep trailing space)
$d.Content.Find.Execute("(sct_crop_image). HAS TO BE DONE MANUALLY", $true, $false, $false, $false, $false, `
    $true, 1, $false, "(sct_crop_image). ", 2) | Out-Null

# 3. Remove ". HAS TO BE DONE MANUALLY" before ": start by generating" in paragraph 2
$d.Content.Find.Execute("fslmaths -add ). HAS TO BE DONE MANUALLY : start by generating", $true, $false, $false, $false, $false, `
    $true, 1, $false, "fslmaths -add ) : start by generating", 2) | Out-Null

# 4. Fix typo "lore" -> "lower"
$d.Content.Find.Execute("the upper and lore nonzero points", $true, $false, $false, $false, $false, `
    $true, 1, $false, "the upper and lower nonzero points", 2) | Out-Null

# 5. Remove the "((Dont do this ...))" paragraph plus the 3 blank paragraphs that follow it
$rStart = $d.Content.Duplicate
$rStart.Find.Execute("((Dont do this", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$rStart.Expand(4) | Out-Null

$rEnd = $d.Content.Duplicate
$rEnd.Find.Execute("8 - Create a mask", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$rEnd.Expand(4) | Out-Null

$d.Range($rStart.Start, $rEnd.Start).Delete() | Out-Null

# 6. Bold "You only have to do this once" inside paragraph 9
$rBold = $d.Content.Duplicate
$rBold.Find.Execute("You only have to do this once", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$rBold.Font.Bold = 1

# 7. Merge the "11 - Crop the straight centerline ..." paragraph with the following
#    "cropped centerline into the template space ..." paragraph into a single paragraph
$r11 = $d.Content.Duplicate
$r11.Find.Execute("11 - Crop the straight centerline", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$r11.Expand(4) | Out-Null
$markRange = $d.Range($r11.End - 1, $r11.End)
$markRange.Delete() | Out-Null

# 8. Move the _GoBack bookmark from the end of the "For T1 volumes..." paragraph to
#    right after "...push this straight " (between the two runs of paragraph 11)
$rAfterStraight = $d.Content.Duplicate
$rAfterStraight.Find.Execute("push this straight ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete() | Out-Null
}
$d.Bookmarks.Add("_GoBack", $d.Range($rAfterStraight.End, $rAfterStraight.End)) | Out-Null

# 9. Append additional sentence to paragraph "12 - use this centerline..."
$r12 = $d.Content.Duplicate
$r12.Find.Execute("12 - use this centerline and the volume to normalize intensity (sct_normalize )", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Range($r12.End, $r12.End).InsertAfter(". Before you should apply the transformation outputed in 10 to the centerline generated in 11") | Out-Null

# 10. Merge the two runs of the final "For T1 volumes..." paragraph into a single run
$rT1 = $d.Content.Duplicate
$rT1.Find.Execute("For T1 volumes", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$rT1.Expand(4) | Out-Null
$t1Text = $rT1.Text
if ($t1Text.EndsWith([char]13)) {
    $t1Text = $t1Text.Substring(0, $t1Text.Length - 1)
}
$rT1.Text = $t1Text

